# Updates the regression-coefficient figures reported in the
# "multiplelinreg" page (scikit-learn, custom OLS, gradient-descent and
# stochastic-gradient-descent sections all get refreshed numbers).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $found = $find.Execute($old, $true, $false, $false, $false, $false, `
                            $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find.Execute could not locate text: '$old'"
    }
}

# scikit-learn fit (appears twice: "Using scikit-learn" and the duplicated
# "Using Custom Library OLS" summary paragraph share identical wording)
Replace-Text "1.0307 and coefficients" "0.9863 and coefficients"
Replace-Text "3.0293, and" "3.0371, and"
Replace-Text "2.02" "1.9549"

# gradient-descent fit
Replace-Text "1.0298 and coefficients" "0.9857 and coefficients"
Replace-Text "3.0289, and" "3.0359, and"
Replace-Text "2.0193" "1.9543"

# stochastic-gradient-descent fit (numpy array repr)
Replace-Text "array([1.01741031]) and coefficients" "array([0.99955242]) and coefficients"
Replace-Text "array([3.03917656]), and" "array([3.0286027]), and"
Replace-Text "array([2.03102484])" "array([1.93519378])"
